# ADD results from server
# Update computed investment-cost results for each year's sheet with the
# latest values returned by the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 8298.121061896434
$ws.Range("E2").Value = 231408.4314442363
$ws.Range("G2").Value = 64767.40570129467
$ws.Range("I2").Value = 129368.5348562943
$ws.Range("L2").Value = 387937.771360128
$ws.Range("M2").Value = 87109.77021824001
$ws.Range("N2").Value = 56002.98791663077
$ws.Range("O2").Value = 55795.91474765583

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 24962.59086493072
$ws.Range("E2").Value = 136526.8409003387
$ws.Range("I2").Value = 167264.4907388069
$ws.Range("L2").Value = 50814.49290518981
$ws.Range("M2").Value = 56347.39187164272
$ws.Range("N2").Value = 15428.15328371197
$ws.Range("O2").Value = 21626.51095245816

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22034.54043650656
$ws.Range("B2").Value = 17690.57314618799
$ws.Range("E2").Value = 91724.35221653034
$ws.Range("I2").Value = 123092.8689171273
$ws.Range("M2").Value = 36699.26852609014
$ws.Range("N2").Value = 31741.5082371194
$ws.Range("O2").Value = 25048.83495982233

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 914.0641520319534
$ws.Range("O2").Value = 0

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 23670.66806629021
$ws.Range("N2").Value = 3478.034812508252
$ws.Range("O2").Value = 16343.81581915547
